$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3179.4
$ws.Range("I53").Value = 115.4
$ws.Range("J53").Value = 6243.4
$ws.Range("K53").Value = 115.4
$ws.Range("L53").Value = 6243.4
$ws.Range("M53").Value = 521.6
$ws.Range("N53").Value = -7517.4

$ws.Range("H98").Value = 780.90625
$ws.Range("I98").Value = 634.96155
$ws.Range("J98").Value = 1413.3334
$ws.Range("K98").Value = 634.96155
$ws.Range("L98").Value = 1413.3334
$ws.Range("M98").Value = 863.03845
$ws.Range("N98").Value = -4409.3334

$ws.Range("H122").Value = 780.90625
$ws.Range("I122").Value = 634.96155
$ws.Range("J122").Value = 1413.3334
$ws.Range("K122").Value = 1904.88465
$ws.Range("L122").Value = 4240.0002
$ws.Range("M122").Value = 545.11535
$ws.Range("N122").Value = -9140.0002

$ws.Range("H129").Value = 157154.48
$ws.Range("J129").Value = 182832.84
$ws.Range("L129").Value = 548498.52
$ws.Range("N129").Value = -558498.52

$ws.Range("H138").Value = 2180.9424
$ws.Range("I138").Value = 1369.65
$ws.Range("J138").Value = 2688
$ws.Range("K138").Value = 4108.950000000001
$ws.Range("L138").Value = 8064
$ws.Range("M138").Value = 1031.049999999999
$ws.Range("N138").Value = -18344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644.84375
$ws.Range("I2").Value = 627.4
$ws.Range("K2").Value = 627.4
$ws.Range("M2").Value = -514.4

$ws.Range("H32").Value = 6898.768
$ws.Range("I32").Value = 5641
$ws.Range("J32").Value = 16489.25
$ws.Range("K32").Value = 5641
$ws.Range("L32").Value = 16489.25
$ws.Range("M32").Value = -5354
$ws.Range("N32").Value = -17063.25

$ws.Range("H61").Value = 2384.1025
$ws.Range("I61").Value = 2205.0303
$ws.Range("J61").Value = 3369
$ws.Range("K61").Value = 2205.0303
$ws.Range("L61").Value = 3369
$ws.Range("M61").Value = -1993.0303
$ws.Range("N61").Value = -3793

$ws.Range("H63").Value = 4466307
$ws.Range("I63").Value = 2358.3333
$ws.Range("K63").Value = 2358.3333
$ws.Range("M63").Value = -1672.3333

$ws.Range("H66").Value = 4466307
$ws.Range("I66").Value = 2358.3333
$ws.Range("K66").Value = 11791.6665
$ws.Range("M66").Value = -8359.666499999999

$ws.Range("H116").Value = 644.84375
$ws.Range("I116").Value = 627.4
$ws.Range("K116").Value = 627.4
$ws.Range("M116").Value = 1666.6

$ws.Range("H132").Value = 10539.821
$ws.Range("I132").Value = 1363.9183
$ws.Range("K132").Value = 4091.7549
$ws.Range("M132").Value = -1561.7549

$ws.Range("H136").Value = 2384.1025
$ws.Range("I136").Value = 2205.0303
$ws.Range("J136").Value = 3369
$ws.Range("K136").Value = 6615.090899999999
$ws.Range("L136").Value = 10107
$ws.Range("M136").Value = -4065.090899999999
$ws.Range("N136").Value = -15207

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644.84375
$ws.Range("I3").Value = 627.4
$ws.Range("K3").Value = 627.4
$ws.Range("M3").Value = -513.4

$ws.Range("H105").Value = 2175765.8
$ws.Range("I105").Value = 1900
$ws.Range("J105").Value = 2501845.5
$ws.Range("K105").Value = 1900
$ws.Range("L105").Value = 2501845.5
$ws.Range("M105").Value = -153
$ws.Range("N105").Value = -2505339.5

$ws.Range("H134").Value = 3672.4524
$ws.Range("I134").Value = 3656.9429
$ws.Range("K134").Value = 10970.8287
$ws.Range("M134").Value = -8435.8287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 17.375
$ws.Range("I7").Value = 12.111111
$ws.Range("J7").Value = 24.142857
$ws.Range("K7").Value = 12.111111
$ws.Range("L7").Value = 24.142857
$ws.Range("M7").Value = 100.888889
$ws.Range("N7").Value = -250.142857

$ws.Range("H105").Value = 1143.9231
$ws.Range("I105").Value = 996.36365
$ws.Range("K105").Value = 996.36365
$ws.Range("M105").Value = 750.63635

$ws.Range("H107").Value = 1222.2693
$ws.Range("I107").Value = 456.33334
$ws.Range("K107").Value = 456.33334
$ws.Range("M107").Value = 1463.66666

$ws.Range("H132").Value = 1653.8596
$ws.Range("I132").Value = 1291.1818
$ws.Range("J132").Value = 2881.3845
$ws.Range("K132").Value = 3873.5454
$ws.Range("L132").Value = 8644.1535
$ws.Range("M132").Value = -1343.5454
$ws.Range("N132").Value = -13704.1535

$ws.Range("H134").Value = 884.1429000000001
$ws.Range("I134").Value = 816.1818
$ws.Range("J134").Value = 1133.3334
$ws.Range("K134").Value = 2448.5454
$ws.Range("L134").Value = 3400.0002
$ws.Range("M134").Value = 86.45460000000003
$ws.Range("N134").Value = -8470.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 366.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 366.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1099.5
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -1323.5

$ws.Range("H131").Value = 684.78
$ws.Range("J131").Value = 711.8461
$ws.Range("L131").Value = 2135.5383
$ws.Range("N131").Value = -12215.5383

$ws.Range("H132").Value = 626.5
$ws.Range("J132").Value = 676.25
$ws.Range("L132").Value = 6086.25
$ws.Range("N132").Value = -11146.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 19000
$ws.Range("J69").Value = 19000
$ws.Range("L69").Value = 19000
$ws.Range("N69").Value = -20498

$ws.Range("H72").Value = 19000
$ws.Range("J72").Value = 19000
$ws.Range("L72").Value = 57000
$ws.Range("N72").Value = -64488

$ws.Range("H126").Value = 2982.151
$ws.Range("I126").Value = 2477.7576
$ws.Range("J126").Value = 3814.4
$ws.Range("K126").Value = 7433.2728
$ws.Range("L126").Value = 11443.2
$ws.Range("M126").Value = -4963.2728
$ws.Range("N126").Value = -16383.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4033.9375
$ws.Range("I7").Value = 3603.3845
$ws.Range("J7").Value = 5899.6665
$ws.Range("K7").Value = 3603.3845
$ws.Range("L7").Value = 5899.6665
$ws.Range("M7").Value = -3491.3845
$ws.Range("N7").Value = -6123.6665

$ws.Range("H82").Value = 1492.3636
$ws.Range("I82").Value = 1304.125
$ws.Range("K82").Value = 1304.125
$ws.Range("M82").Value = -943.125

$ws.Range("H85").Value = 1492.3636
$ws.Range("I85").Value = 1304.125
$ws.Range("K85").Value = 1304.125
$ws.Range("M85").Value = -56.125

$ws.Range("H126").Value = 4033.9375
$ws.Range("I126").Value = 3603.3845
$ws.Range("J126").Value = 5899.6665
$ws.Range("K126").Value = 10810.1535
$ws.Range("L126").Value = 17698.9995
$ws.Range("M126").Value = -8340.1535
$ws.Range("N126").Value = -22638.9995

$ws.Range("H132").Value = 390614.88
$ws.Range("I132").Value = 574923.5
$ws.Range("J132").Value = 3566.8
$ws.Range("K132").Value = 1724770.5
$ws.Range("L132").Value = 10700.4
$ws.Range("M132").Value = -1722240.5
$ws.Range("N132").Value = -15760.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1530.4375
$ws.Range("I126").Value = 1128.6666
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 3385.9998
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -915.9998000000001
$ws.Range("N126").Value = -16040
